# Add Table S6 - HostVsymb PERMANOVA worksheet (host vs symbiont PCA PERMANOVA results)
$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Table S6 - HostVsymb PERMANOVA"

# Column widths (characters), matching the layout used by the other PERMANOVA tables
$ws.Columns.Item(1).ColumnWidth = 15.876666666666667
$ws.Columns.Item(2).ColumnWidth = 2.8766666666666665
$ws.Columns.Item(3).ColumnWidth = 13.876666666666667
$ws.Columns.Item(4).ColumnWidth = 4.876666666666666
$ws.Columns.Item(5).ColumnWidth = 4.876666666666666
$ws.Columns.Item(6).ColumnWidth = 6.876666666666666
$ws.Columns.Item(7).ColumnWidth = 2.8766666666666665
$ws.Columns.Item(8).ColumnWidth = 14.876666666666667
$ws.Columns.Item(9).ColumnWidth = 4.876666666666666
$ws.Columns.Item(10).ColumnWidth = 4.876666666666666
$ws.Columns.Item(11).ColumnWidth = 7.876666666666668
$ws.Columns.Item(12).ColumnWidth = 6.876666666666666

# Header row + data rows for the two side-by-side PERMANOVA blocks (host vs symbiont)
$ws.Range("A1").Value = " "
$ws.Range("B1").Value = "Df"
$ws.Range("C1").Value = "Sum of Squares"
$ws.Range("D1").Value = "R2"
$ws.Range("E1").Value = "F"
$ws.Range("F1").Value = "P-value"
$ws.Range("G1").Value = "Df "
$ws.Range("H1").Value = "Sum of Squares "
$ws.Range("I1").Value = "R2 "
$ws.Range("J1").Value = "F "
$ws.Range("K1").Value = "P-value "
$ws.Range("L1").Value = "species"
$ws.Range("A2").Value = "pCO2"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.019
$ws.Range("E2").Value = 0.55
$ws.Range("F2").Value = 0.74883
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 150227
$ws.Range("I2").Value = 0.272
$ws.Range("J2").Value = 11.45
$ws.Range("K2").Value = 0.00067
$ws.Range("L2").Value = "SSID"
$ws.Range("A3").Value = "temperature"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0.075
$ws.Range("E3").Value = 6.65
$ws.Range("F3").Value = 0.004
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 24933
$ws.Range("I3").Value = 0.045
$ws.Range("J3").Value = 5.7
$ws.Range("K3").Value = 0.004
$ws.Range("L3").Value = "SSID"
$ws.Range("A4").Value = "reef environment"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.006
$ws.Range("E4").Value = 0.52
$ws.Range("F4").Value = 0.58294
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 26680
$ws.Range("I4").Value = 0.048
$ws.Range("J4").Value = 6.1
$ws.Range("K4").Value = 0.004
$ws.Range("L4").Value = "SSID"
$ws.Range("A5").Value = "Residual"
$ws.Range("B5").Value = 80
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 0.901
$ws.Range("G5").Value = 80
$ws.Range("H5").Value = 349805
$ws.Range("I5").Value = 0.634
$ws.Range("L5").Value = "SSID"
$ws.Range("A6").Value = "Total"
$ws.Range("B6").Value = 85
$ws.Range("C6").Value = 34
$ws.Range("D6").Value = 1
$ws.Range("G6").Value = 85
$ws.Range("H6").Value = 551646
$ws.Range("I6").Value = 1
$ws.Range("L6").Value = "SSID"
$ws.Range("A7").Value = "pCO2"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0.041
$ws.Range("E7").Value = 1.23
$ws.Range("F7").Value = 0.28714
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 196108
$ws.Range("I7").Value = 0.096
$ws.Range("J7").Value = 4.41
$ws.Range("K7").Value = 0.00333
$ws.Range("L7").Value = "PSTR"
$ws.Range("A8").Value = "temperature"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 0.147
$ws.Range("E8").Value = 13.12
$ws.Range("F8").Value = 0.00067
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 619660
$ws.Range("I8").Value = 0.304
$ws.Range("J8").Value = 41.82
$ws.Range("K8").Value = 0.00067
$ws.Range("L8").Value = "PSTR"
$ws.Range("A9").Value = "reef environment"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.02
$ws.Range("E9").Value = 1.75
$ws.Range("F9").Value = 0.15923
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 167947
$ws.Range("I9").Value = 0.083
$ws.Range("J9").Value = 11.34
$ws.Range("K9").Value = 0.00133
$ws.Range("L9").Value = "PSTR"
$ws.Range("A10").Value = "Residual"
$ws.Range("B10").Value = 71
$ws.Range("C10").Value = 14
$ws.Range("D10").Value = 0.793
$ws.Range("G10").Value = 71
$ws.Range("H10").Value = 1051979
$ws.Range("I10").Value = 0.517
$ws.Range("L10").Value = "PSTR"
$ws.Range("A11").Value = "Total"
$ws.Range("B11").Value = 76
$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 1
$ws.Range("G11").Value = 76
$ws.Range("H11").Value = 2035695
$ws.Range("I11").Value = 1
$ws.Range("L11").Value = "PSTR"
$ws.Range("A12").Value = "pCO2"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 0.136
$ws.Range("E12").Value = 3.48
$ws.Range("F12").Value = 0.01865
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 100597
$ws.Range("I12").Value = 0.298
$ws.Range("J12").Value = 11.23
$ws.Range("K12").Value = 0.00067
$ws.Range("L12").Value = "PAST"
$ws.Range("A13").Value = "temperature"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0.036
$ws.Range("E13").Value = 2.76
$ws.Range("F13").Value = 0.1006
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 51174
$ws.Range("I13").Value = 0.151
$ws.Range("J13").Value = 17.13
$ws.Range("K13").Value = 0.00067
$ws.Range("L13").Value = "PAST"
$ws.Range("A14").Value = "reef environment"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0.021
$ws.Range("E14").Value = 1.64
$ws.Range("F14").Value = 0.20053
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 969
$ws.Range("I14").Value = 0.003
$ws.Range("J14").Value = 0.32
$ws.Range("K14").Value = 0.68754
$ws.Range("L14").Value = "PAST"
$ws.Range("A15").Value = "Residual"
$ws.Range("B15").Value = 62
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 0.807
$ws.Range("G15").Value = 62
$ws.Range("H15").Value = 185165
$ws.Range("I15").Value = 0.548
$ws.Range("L15").Value = "PAST"
$ws.Range("A16").Value = "Total"
$ws.Range("B16").Value = 67
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 1
$ws.Range("G16").Value = 67
$ws.Range("H16").Value = 337905
$ws.Range("I16").Value = 1
$ws.Range("L16").Value = "PAST"

# Residual/Total rows have no F / P-value -- leave those cells explicitly blank
foreach ($addr in @("E5","F5","J5","K5","E6","F6","J6","K6","E10","F10","J10","K10","E11","F11","J11","K11","E15","F15","J15","K15","E16","F16","J16","K16")) {
    $ws.Range($addr).Font.Bold = $false
}

# Match page setup used by the other tables in the workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the original active sheet/selection state
$wb.Worksheets.Item(1).Activate()
